# Generate Report for Handoff
# Updates localization status from "In Translation" to "Ready for handoff"
# and refreshes the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# timestamps on all three sheets, widening the Status columns to fit the
# new, longer status text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-22 06:38:35"
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-22 06:38:31"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3333333333333

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-22 06:38:35"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3333333333333

Write-Output "Report regenerated for handoff"
